# Updated cryptos list - applies latest Price / Volume(1h) figures, and
# corrects the ordering of the OKB / Bittensor rows (47-48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 47/48: Bittensor and OKB swapped places ---
# Row 47 was OKB, becomes Bittensor; row 48 was Bittensor, becomes OKB.
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"

# --- Column D (Price) updates ---
# Leading apostrophe forces text storage so values like "587.01" or "1.00"
# are not reinterpreted as numbers, matching the source's text formatting.
$ws.Range("D2").Value = "'66.969.63"
$ws.Range("D3").Value = "'3.452.45"
$ws.Range("D5").Value = "'587.01"
$ws.Range("D6").Value = "'179.71"
$ws.Range("D7").Value = "'0.631"
$ws.Range("D9").Value = "'3.451.64"
$ws.Range("D11").Value = "'6.98"
$ws.Range("D13").Value = "'4.054.20"
$ws.Range("D15").Value = "'30.19"
$ws.Range("D16").Value = "'66.926.02"
$ws.Range("D18").Value = "'3.453.11"
$ws.Range("D19").Value = "'5.98"
$ws.Range("D20").Value = "'13.94"
$ws.Range("D21").Value = "'374.18"
$ws.Range("D22").Value = "'7.71"
$ws.Range("D23").Value = "'73.70"
$ws.Range("D26").Value = "'0.538"
$ws.Range("D27").Value = "'10.03"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D30").Value = "'5.92"
$ws.Range("D32").Value = "'23.77"
$ws.Range("D34").Value = "'1.30"
$ws.Range("D35").Value = "'7.15"
$ws.Range("D37").Value = "'162.98"
$ws.Range("D38").Value = "'0.884"
$ws.Range("D39").Value = "'28.04"
$ws.Range("D40").Value = "'1.82"
$ws.Range("D41").Value = "'2.67"
$ws.Range("D43").Value = "'2.762.18"
$ws.Range("D44").Value = "'6.39"
$ws.Range("D45").Value = "'0.0702"
$ws.Range("D46").Value = "'25.78"
$ws.Range("D47").Value = "'340.07"
$ws.Range("D48").Value = "'40.25"
$ws.Range("D51").Value = "'32.09"

# --- Column E (Volume (1h) %) updates ---
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("E7").Value = "  +4.93%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("E21").Value = "  -2.54%  "
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("E24").Value = "  +7.45%  "
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("E32").Value = "  -2.66%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("E37").Value = "  +1.37%  "
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("E39").Value = "  -6.57%  "
$ws.Range("E40").Value = "  +0.93%  "
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  +3.97%  "
$ws.Range("E47").Value = "  +6.61%  "
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("E50").Value = "  +2.49%  "
$ws.Range("E51").Value = "  +2.82%  "

Write-Output "cryptos list updated"
